$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire "M" column block (was the "AYED RED" label at M18,
# plus the blank formatted cells below it) - shifts the used range back
# down to column K.
$ws.Range("M18:M21").Delete()

# Clear the remaining "AYED RED" schedule entries (and the TEORIA/PRACTICA/
# room cells that went with that block), reverting each cell back to the
# plain bordered look used throughout the rest of the grid.
$formatDonor = $ws.Range("C18")
$cellsToClear = @("B19", "B20", "B21", "D21", "D22", "D23", "E24", "E25", "E26")

foreach ($addr in $cellsToClear) {
    $cell = $ws.Range($addr)
    $cell.ClearContents()
    $formatDonor.Copy()
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update the active selection, matching the saved view state.
$ws.Range("M27").Select()
